# Insert a new weekly data row at row 15, pushing all subsequent rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 15..131 down to 16..132 by inserting a blank row at 15.
$ws.Rows("15:15").Insert()

# Populate the newly inserted row 15 with the new weekly observation.
$ws.Cells.Item(15, 1).Value = 9
$ws.Cells.Item(15, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(15, 3).Value = "Metropolitana"
$ws.Cells.Item(15, 4).Value = 44530
$ws.Cells.Item(15, 5).Value = 13
$ws.Cells.Item(15, 6).Value = "Fruta"
$ws.Cells.Item(15, 7).Value = 100101
$ws.Cells.Item(15, 8).Value = "Berries"
$ws.Cells.Item(15, 9).Value = 100101001
$ws.Cells.Item(15, 10).Value = "Arándano (blue)"
$ws.Cells.Item(15, 11).Value = "Sin especificar"
$ws.Cells.Item(15, 12).Value = "Primera"
$ws.Cells.Item(15, 13).Value = 550
$ws.Cells.Item(15, 14).Value = 4800
$ws.Cells.Item(15, 15).Value = 5000
$ws.Cells.Item(15, 16).Value = 4927
$ws.Cells.Item(15, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(15, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(15, 19).Value = 2464
$ws.Cells.Item(15, 20).Value = 2
